$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.412.99"
$ws.Cells.Item(2, 5).Value = "  +0.54%  "
$ws.Cells.Item(3, 4).Value = "1.611.04"
$ws.Cells.Item(3, 5).Value = "  +1.08%  "
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.18"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.41%  "
$ws.Cells.Item(6, 5).Value = "  -0.44%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$ws.Cells.Item(8, 5).Value = "  -0.34%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.23"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.26%  "
$ws.Cells.Item(11, 5).Value = "  -0.42%  "
$ws.Cells.Item(12, 4).Value = "1.837.88"
$ws.Cells.Item(12, 5).Value = "  +1.14%  "
$ws.Cells.Item(13, 4).Value = "1.612.90"
$ws.Cells.Item(13, 5).Value = "  +2.09%  "
$ws.Cells.Item(14, 5).Value = "  -0.04%  "
$ws.Cells.Item(15, 5).Value = "  +0.01%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "63.58"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.39%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "235.59"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +9.45%  "
$ws.Cells.Item(18, 4).Value = "26.414.34"
$ws.Cells.Item(18, 5).Value = "  +0.60%  "
$ws.Cells.Item(19, 5).Value = "  +0.54%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.66"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +4.29%  "
$ws.Cells.Item(21, 5).Value = "  -0.11%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.27"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.49%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.19"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +3.85%  "
$ws.Cells.Item(24, 5).Value = "  +0.17%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "146.77"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.03%  "
$ws.Cells.Item(27, 5).Value = "  +0.36%  "
$ws.Cells.Item(28, 5).Value = "  +0.06%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.49"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.54%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0495"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.20%  "
$ws.Cells.Item(31, 5).Value = "  -0.73%  "
$ws.Cells.Item(32, 4).Value = "1.506.98"
$ws.Cells.Item(32, 5).Value = "  +6.46%  "
$ws.Cells.Item(33, 5).Value = "  +1.23%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.94"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.95%  "
$ws.Cells.Item(36, 5).Value = "  -0.32%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.562"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.66%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0165"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.03%  "
$ws.Cells.Item(40, 5).Value = "  +0.96%  "
$ws.Cells.Item(41, 5).Value = "  -0.08%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.19"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.12%  "
$ws.Cells.Item(43, 4).Value = "1.750.50"
$ws.Cells.Item(43, 5).Value = "  +1.24%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.761"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.05%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.918"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -2.54%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "61.33"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.72%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "89.94"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.86%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0103"
$ws.Cells.Item(48, 5).Value = "  -2.34%  "
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.49"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.52%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0501"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.14%  "
$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0959"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.59%  "
